# "Actualización desde MV -datos-"
# Updates a handful of previously-published index values (small revisions)
# and appends the newly published month (01-08-2021) as a new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised values in column B (Ventas SOFOFA base 2014=100) -------------
$ws.Cells.Item(182, 2).Value = 95.45
$ws.Cells.Item(186, 2).Value = 99.08
$ws.Cells.Item(187, 2).Value = 98.7
$ws.Cells.Item(188, 2).Value = 99.14
$ws.Cells.Item(222, 2).Value = 102.62
$ws.Cells.Item(236, 2).Value = 103.11
$ws.Cells.Item(237, 2).Value = 105.61
$ws.Cells.Item(239, 2).Value = 96.73
$ws.Cells.Item(240, 2).Value = 106.59
$ws.Cells.Item(241, 2).Value = 107.4
$ws.Cells.Item(246, 2).Value = 95.73
$ws.Cells.Item(260, 2).Value = 112.63

# --- New row for the newly published period (01-08-2021) ------------------
$dateCell = $ws.Cells.Item(261, 1)
# Force text storage so "01-08-2021" is kept as a literal string (matching
# the rest of column A) instead of being auto-converted into a date serial.
$dateCell.NumberFormat = "@"
$dateCell.Value = "01-08-2021"
# Restore the default (unstyled) look so the new cell matches its neighbours
# (no explicit style index), same as every other cell in column A.
$dateCell.Style = "Normal"

$ws.Cells.Item(261, 2).Value = 116.32
